$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-24 Monday" "2025-02-25 Tuesday"

Replace-Text "605×7=4235" "420×5=2100"
Replace-Text "358×8=2864" "126×8=1008"
Replace-Text "408×3=1224" "741×7=5187"
Replace-Text "150×5=750" "788×8=6304"
Replace-Text "232×8=1856" "722×3=2166"

Replace-Text "625×8=5000" "675×3=2025"
Replace-Text "539×5=2695" "537×5=2685"
Replace-Text "676×3=2028" "572×9=5148"
Replace-Text "811×8=6488" "330×6=1980"
Replace-Text "715×8=5720" "207×5=1035"

Replace-Text "340×6=2040" "141×7=987"
Replace-Text "902×9=8118" "237×9=2133"
Replace-Text "722×2=1444" "832×5=4160"
Replace-Text "521×4=2084" "326×9=2934"
Replace-Text "903×5=4515" "779×7=5453"

Replace-Text "137×2=274" "223×3=669"
Replace-Text "148×2=296" "511×6=3066"
Replace-Text "729×6=4374" "691×3=2073"
Replace-Text "748×9=6732" "545×9=4905"
Replace-Text "232×3=696" "394×8=3152"

Replace-Text "449×7=3143" "181×3=543"
Replace-Text "162×4=648" "585×5=2925"
Replace-Text "722×6=4332" "318×8=2544"
Replace-Text "803×9=7227" "854×3=2562"
Replace-Text "238×6=1428" "482×9=4338"
